# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new F value for the "展览" sheet
$exhibitUpdates = @{
    3  = 22
    6  = 20
    8  = 77
    9  = 446
    12 = 553
    13 = 25
    14 = 287
    15 = 23
    16 = 347
    21 = 88
    22 = 863
    23 = 1378
    25 = 313
    27 = 65
    30 = 83
    31 = 206
    32 = 238
    34 = 1598
    37 = 152
    38 = 569
    39 = 298
    40 = 3459
    41 = 409
    42 = 181
    43 = 885
    45 = 59
}

# Row -> new F value for the "全部类型" sheet
# (identical to the above, except row 40 which ends up one higher)
$allTypesUpdates = @{
    3  = 22
    6  = 20
    8  = 77
    9  = 446
    12 = 553
    13 = 25
    14 = 287
    15 = 23
    16 = 347
    21 = 88
    22 = 863
    23 = 1378
    25 = 313
    27 = 65
    30 = 83
    31 = 206
    32 = 238
    34 = 1598
    37 = 152
    38 = 569
    39 = 298
    40 = 3460
    41 = 409
    42 = 181
    43 = 885
    45 = 59
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}

Write-Output "Updated $($exhibitUpdates.Count) rows on 展览 and $($allTypesUpdates.Count) rows on 全部类型"
